$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.748.31'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '3.100.61'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.58'
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.87'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.099.15'
$ws.Range("E8").Value = '  +1.12%  '
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.22'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("D13").Value = '3.633.79'
$ws.Range("E13").Value = '  +1.27%  '
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.74'
$ws.Range("E15").Value = '  -2.38%  '
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("D17").Value = '57.783.79'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").Value = '3.099.22'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.09'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.03'
$ws.Range("E21").Value = '  -1.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '340.44'
$ws.Range("E22").Value = '  +2.65%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.85'
$ws.Range("E25").Value = '  +2.32%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '0.0₃0919'
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.19'
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.98'
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("E34").Value = '  -0.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.08'
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.63'
$ws.Range("E36").Value = '  +2.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.14'
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.10'
$ws.Range("E38").Value = '  -1.56%  '
$ws.Range("E39").Value = '  -1.19%  '
$ws.Range("E40").Value = '  -2.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.97'
$ws.Range("E41").Value = '  +1.57%  '
$ws.Range("D42").Value = '3.140.53'
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("E43").Value = '  +3.80%  '
$ws.Range("E44").Value = '  +9.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.79'
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").Value = '2.294.61'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("E49").Value = '  +4.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.54'
$ws.Range("E50").Value = '  -1.48%  '
$ws.Range("E51").Value = '  +1.56%  '
